$wb = $excel.ActiveWorkbook

# --- "Reguły" sheet (rule descriptions) — reorder the object lists in the
# rule text for rules 1, 2, 3, 5, 6, 7 (rule 4 / B5 is unchanged) ---
$wsReguly = $wb.Worksheets.Item(8)
$wsReguly.Range("B2").Value = "(attempts >=  3.0) & (pregnancy <=  0.0) => (class <= 1) ['a3', 'a7', 'a1']"
$wsReguly.Range("B3").Value = "(sperm >=  3.0) => (class <= 1) ['a22', 'a25']"
$wsReguly.Range("B4").Value = "(age >=  40.0) & (pregnancy <=  0.0) => (class <= 1) ['a3', 'a15']"
$wsReguly.Range("B6").Value = "(age >=  42.0) => (class <= 1) ['a14', 'a3']"
$wsReguly.Range("B7").Value = "(age <=  31.0) & (attempts <=  1.0) & (endometrium <=  1.0) => (class >= 2) ['a24', 'a9', 'a11', 'a12']"
$wsReguly.Range("B8").Value = "(frozen_embryos >=  8.0) & (sperm <=  1.0) => (class >= 2) ['a16', 'a6']"

# --- "Statystyki reguł" sheet — updated coverage values for rule 1 and rule 3 ---
$wsStatystyki = $wb.Worksheets.Item(9)
$wsStatystyki.Range("C2").Value = 0.375
$wsStatystyki.Range("C4").Value = 0.25

# --- "Walidacja krzyżowa" sheet — reorder the metric rows (values stay
# attached to their original labels, only the row order changes) ---
$wsWalidacja = $wb.Worksheets.Item(10)
$wsWalidacja.Range("A1").Value = "accuracy"
$wsWalidacja.Range("B1").Value = 0.36
$wsWalidacja.Range("A2").Value = "not_classified"
$wsWalidacja.Range("B2").Value = 0.56
$wsWalidacja.Range("A3").Value = "correct"
$wsWalidacja.Range("B3").Value = 0.8181818181818182
$wsWalidacja.Range("A4").Value = "f1_score"
$wsWalidacja.Range("B4").Value = 0.48
